$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.336.04'
$ws.Range('E2').Value = '  -0.09%  '
$ws.Range('D3').Value = '3.140.84'
$ws.Range('E3').Value = '  -0.43%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '610.74'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.77%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.95'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.30%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').Value = '3.138.49'
$ws.Range('E8').Value = '  -0.42%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.529'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.32%  '
$ws.Range('E10').Value = '  +0.14%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.39'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.64%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.477'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.37%  '
$ws.Range('E13').Value = '  +2.54%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.67'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.53%  '
$ws.Range('D15').Value = '3.656.94'
$ws.Range('E15').Value = '  -0.43%  '
$ws.Range('E16').Value = '  +2.46%  '
$ws.Range('D17').Value = '64.341.81'
$ws.Range('D18').Value = '3.139.25'
$ws.Range('E18').Value = '  -0.41%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.88'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.18%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '477.69'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.45%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.73'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.39%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.724'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.98%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.83'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.27%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.69'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.54%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '85.25'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E26').Value = '  -0.01%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.62'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.11%  '
$ws.Range('E28').Value = '  -2.90%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.47'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +9.03%  '
$ws.Range('E30').Value = '  +2.84%  '
$ws.Range('E31').Value = '  -4.81%  '
$ws.Range('E32').Value = '  -0.13%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '26.69'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.66%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.64'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.71%  '
$ws.Range('E35').Value = '  +0.43%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.96'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.70%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '52.48'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').Value = '0.0₃0747'
$ws.Range('E38').Value = '  +4.59%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '456.43'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.59%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.03'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +4.85%  '
$ws.Range('E41').Value = '  +0.06%  '
$ws.Range('E42').Value = '  -0.49%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.35'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.26%  '
$ws.Range('D44').Value = '2.865.78'
$ws.Range('E44').Value = '  +0.86%  '
$ws.Range('E45').Value = '  -0.65%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.28'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.14%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.44'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +4.98%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '26.57'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.42%  '
$ws.Range('E49').Value = '  +0.08%  '
$ws.Range('E50').Value = '  -0.03%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '120.56'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.14%  '
